$d = $word.ActiveDocument

function Replace-ParagraphXml($anchorText, $newInnerXml, $paraAttrs) {
    $r = $d.Content
    $found = $r.Find.Execute($anchorText)
    if (-not $found) {
        throw "Anchor not found: $anchorText"
    }
    $para = $r.Paragraphs(1)
    $pr = $para.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p ' + $paraAttrs + '>' + $newInnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pr.InsertXML($xml)
}

# --- 1. "Guiseppe: " -> spell-checked "Guiseppe" + ": " (first occurrence, task table) ---
$attrs1 = 'w14:paraId="134AF79C" w14:textId="588E6191" w:rsidR="00782380" w:rsidRDefault="00667F7F" w:rsidP="00F44F3A"'
$inner1 = '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Guiseppe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="008640AB"><w:t>Begin Analysis document</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Andrew: </w:t></w:r><w:r w:rsidR="008640AB"><w:t>Revise database schemas</w:t></w:r>'
Replace-ParagraphXml "Guiseppe: " $inner1 $attrs1

# --- 2. "Complete search filters..." -> "Begin creating" + " search filters..." ---
$attrs2 = 'w14:paraId="7863403D" w14:textId="164F7BC8" w:rsidR="00782380" w:rsidRDefault="008640AB" w:rsidP="00F44F3A"'
$inner2 = '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr><w:r><w:t>Begin creating</w:t></w:r><w:r><w:t xml:space="preserve"> search filters on all screens and review all completed wireframes</w:t></w:r>'
Replace-ParagraphXml "Complete search filters on all screens and review all completed wireframes" $inner2 $attrs2

# --- 3. "Create base wireframes..." -> "Begin creating" + " base wireframes..." ---
$attrs3 = 'w14:paraId="6241F3B4" w14:textId="021C2541" w:rsidR="00782380" w:rsidRDefault="008640AB" w:rsidP="00F44F3A"'
$inner3 = '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr><w:r><w:t>Begin creating</w:t></w:r><w:r><w:t xml:space="preserve"> base wireframes for all major navigation screens</w:t></w:r>'
Replace-ParagraphXml "Create base wireframes for all major navigation screens" $inner3 $attrs3

# --- 4. "Create wireframes for team screen and complete technical design document" -> multi-run split with relocated _GoBack bookmark ---
$attrs4 = 'w14:paraId="6E7C8A6D" w14:textId="1B61C4FD" w:rsidR="00782380" w:rsidRDefault="008640AB" w:rsidP="00F44F3A"'
$inner4 = '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr><w:r><w:t>Begin c</w:t></w:r><w:bookmarkStart w:id="5" w:name="_GoBack"/><w:bookmarkEnd w:id="5"/><w:r><w:t>reate wireframes for team screen</w:t></w:r><w:r><w:t xml:space="preserve"> and ready screen</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:t>begin</w:t></w:r><w:r><w:t xml:space="preserve"> technical design document</w:t></w:r>'
Replace-ParagraphXml "Create wireframes for team screen and complete technical design document" $inner4 $attrs4

# --- 5. "Guiseppe Ragusa" -> spell-checked "Guiseppe" + " Ragusa" ---
$attrs5 = 'w14:paraId="3AD77AA6" w14:textId="2575765E" w:rsidR="00782380" w:rsidRDefault="00325770" w:rsidP="00F44F3A"'
$inner5 = '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Guiseppe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Ragusa</w:t></w:r>'
Replace-ParagraphXml "Guiseppe Ragusa" $inner5 $attrs5

# --- 6. Remove the _GoBack bookmark from the "Agenda: Revising Wireframes" paragraph (it moved to #4) ---
$attrs6 = 'w14:paraId="0D07477C" w14:textId="6EB5F99C" w:rsidR="00782380" w:rsidRDefault="00782380" w:rsidP="00F44F3A"'
$inner6 = '<w:pPr><w:pStyle w:val="Standard1"/><w:spacing w:before="120" w:after="120"/></w:pPr><w:r><w:t xml:space="preserve">Agenda: </w:t></w:r><w:r w:rsidR="00271C86"><w:t>Revising Wireframes</w:t></w:r>'
Replace-ParagraphXml "Revising Wireframes" $inner6 $attrs6

# --- 7. " Giuseppe.Ragusa@georgebrown" -> " " + spell-checked "Giuseppe.Ragusa@georgebrown" ---
$attrs7 = 'w14:paraId="4C905EC5" w14:textId="08D35F3C" w:rsidR="00782380" w:rsidRPr="0016291E" w:rsidRDefault="00782380" w:rsidP="00782380"'
$inner7 = '<w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="009D6683"><w:rPr><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-GB"/></w:rPr><w:t>Team member 3</w:t></w:r><w:r w:rsidR="002E154B" w:rsidRPr="009D6683"><w:rPr><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="002E154B" w:rsidRPr="009D6683"><w:rPr><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-GB"/></w:rPr><w:t>Giuseppe.Ragusa@georgebrown</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Replace-ParagraphXml "Giuseppe.Ragusa@georgebrown" $inner7 $attrs7

Write-Host "All steps done"
